# Trabalho1-PC2.docx edits
# 1) Delivery date changes from 25/11 to 05/05.
# 2) The group-size paragraph is reworded: "em 7 grupos de 4 pessoas (28 alunos no
#    total)." becomes "em 10 grupos de 3 pessoas.", and "em" is no longer bold.
# 3) "obrigatoriamente " is inserted before "na linguagem Java."
# 4) The (hidden) "_GoBack" bookmark is relocated from after "...em um" to a
#    collapsed position after "...o mesmo p" (inside "para").
# 5) "seus números de matrícula" becomes "seus respectivos números de matrícula".

$d = $word.ActiveDocument

# --- 1) Delivery date -------------------------------------------------------
$d.Content.Find.Execute("25/11", $true, $false, $false, $false, $false, $true, 1, $false, "05/05", 2) | Out-Null

# --- 2) Group size sentence --------------------------------------------------
$d.Content.Find.Execute("em 7 grupos de 4 pessoas (28 alunos no total). ", $true, $false, $false, $false, $false, $true, 1, $false, "em 10 grupos de 3 pessoas. ", 2) | Out-Null

# The replacement inherited the (bold) formatting of "em " from the matched
# text; "em" itself should no longer be bold in the final version, while
# " 10 grupos de 3 pessoas. " stays bold.
$rEm = $d.Content
$rEm.Find.Execute("em 10", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rEmOnly = $d.Range($rEm.Start, $rEm.Start + 2)
$rEmOnly.Font.Bold = 0

# --- 3) "obrigatoriamente" ----------------------------------------------------
$d.Content.Find.Execute("trabalho deverá ser feito n", $true, $false, $false, $false, $false, $true, 1, $false, "trabalho deverá ser feito obrigatoriamente n", 2) | Out-Null

# --- 4) Relocate the "_GoBack" bookmark --------------------------------------
$rAnchor = $d.Content
$rAnchor.Find.Execute("o mesmo p", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($rAnchor.End, $rAnchor.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 5) "respectivos" + dedicated run for "números de matrícula" ------------
$d.Content.Find.Execute("seus números de matrícula", $true, $false, $false, $false, $false, $true, 1, $false, "seus respectivos números de matrícula", 2) | Out-Null

$rNum = $d.Content
$rNum.Find.Execute("números de matrícula", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rNum.Font.Bold = 1
$rNum.Font.Italic = 1

Write-Host "edit.ps1 completed successfully"
